# ConfigTracker2.xlsx update:
# Fix the typo "insall" -> "install" in the fireDomain.install.expense
# item name, which appears in column D of the "Tracker" sheet for the
# two rows that reference the expense variant (rows 2 and 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracker")

$ws.Range("D2").Value = "fireDomain.install.expense"
$ws.Range("D3").Value = "fireDomain.install.expense"
